# "Generate Report for handoff": refresh the localization-status report.
# The row for 27b89366-4817-4392-9fec-e9a51b422d94.md moves from
# "Ready for handoff" to "In Translation" on the Overview sheet and on
# both per-language sheets, and the still-pending "Ready for handoff"
# rows get a fresh "Latest Handoff Datetime" stamp from the new handoff run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B4").Value = "In Translation"
$ov.Range("C4").Value = "In Translation"

# --- zh-cn sheet ---
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("B4").Value = "In Translation"
$zh.Range("D2").Value = "2016-01-26 10:20:35"
$zh.Range("D3").Value = "2016-01-26 10:20:35"
$zh.Range("D5").Value = "2016-01-26 10:20:35"
$zh.Range("D6").Value = "2016-01-26 10:20:35"
$zh.Range("D7").Value = "2016-01-26 10:20:35"

# --- de-de sheet ---
$de = $wb.Worksheets.Item("de-de")
$de.Range("B4").Value = "In Translation"
$de.Range("D2").Value = "2016-01-26 10:20:47"
$de.Range("D3").Value = "2016-01-26 10:20:47"
$de.Range("D5").Value = "2016-01-26 10:20:47"
$de.Range("D6").Value = "2016-01-26 10:20:47"
$de.Range("D7").Value = "2016-01-26 10:20:47"
